# Review_186.docx -> Review_185.docx style edit
# Applies the paper-review content swap described by the commit diff:
#  - new title / paper link / secondary link
#  - four new paragraphs inserted after the secondary link (3 blank + 1 with
#    the new "first paper in the series" intro paragraph)
#  - two body paragraphs get their text replaced
#  - a long run of old body paragraphs is deleted
#  - the final summary paragraph gets new text
#
# Text replacements are done with a paragraph-scoped Find.Execute (rather
# than Range.Text = ...) so the freshly written runs don't inherit a stale
# xml:space="preserve" flag from the run being replaced.

$d = $word.ActiveDocument

function Replace-ParagraphText($paragraph, [string]$oldText, [string]$newText) {
    $rng = $paragraph.Range
    $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
}

# --- 1) Title (Heading1) ----------------------------------------------------
Replace-ParagraphText $d.Paragraphs.Item(1) `
    'Review 186: HiPPO: Recurrent Memory with Optimal Polynomial Projections' `
    'Review 185: Legendre Memory Units: Continuous-Time Representation in Recurrent Neural Networks'

# --- 2) Bold "Paper:" line ---------------------------------------------------
Replace-ParagraphText $d.Paragraphs.Item(2) `
    'Paper: https://arxiv.org/abs/2008.07669v2' `
    'Paper: https://arxiv.org/abs/2102.11417v2'

# --- 3) Paragraph 3 is an empty <w:p/> and is left untouched ----------------

# --- 4) Secondary link line --------------------------------------------------
Replace-ParagraphText $d.Paragraphs.Item(4) `
    'https://arxiv.org/abs/2008.07669' `
    'https://proceedings.neurips.cc/paper_files/paper/2019/file/952285b9b7e7a1be5aa7849f32ffff05-Paper.pdf'

# --- 5) Replace the two surviving Hebrew body paragraphs (before inserting
#        new paragraphs, so their indices -- 6 and 8 -- stay correct) -------
Replace-ParagraphText $d.Paragraphs.Item(6) `
    'הגענו למאמר השני בסדרה - המאמר הזה חשוב מאוד כי הוא מפתח בסיס מתמטי מוצק המשמש כל המודלים מבוססים על מערכות דינמיות לינאריות כולל כמובן ממבה. המאמר הזה קצת (די הרבה) כבד מתמטית אך אנסה לעשות כמיטב יכולתי כדי להעביר לכם את המסר העיקרי שהוא מביא איתו. ' `
    'המאמר בונה מערכת דינמית המתוארת על ידי משוואה דיפרנציאלית לינארית (מערכת דינמית, משוואה 1 במאמר) כאשר (m(t הוא וקטור הזכרון ו- (u(t כאמור הקלט (כרגע חד ממדי). מתברר שעבור בחירה מסוימת של מטריצת A במשוואה של המערכת הדינמית ניתן לתאר את הקלט (בפרק זמן מסוים) על ידי שילוב של פונקציית הזיכרון (m(t ופונקציות מתמטיות הנקראות פולינומים של Legendre (משוואה 3 במאמר). כלומר ניתן לתאר את כל מה שקרה מבחינת הקלט עד זמן מסוים על ידי פונקציה (m(t - וזה בדיוק מה שרצינו, נכון?'

Replace-ParagraphText $d.Paragraphs.Item(8) `
    'בסקירה הקודמת דיברנו על איך ניתן לבנות וקטור זיכרון (m(t בעל יכולת לשחזר פונקצית קלט (u(x ל-x מאינטרוול ; כאן t מסמן גודל חלון הקשר (כלומר אורך הזיכרון). פונקצית (m(t ממודלת על ידי מערכת דינמית לינארית ושילובה עם פולינומי Legendre משחזר לנו את הקלט u. נעיר שאנו עובדים עם הגרסאות הדיסקרטיות של המודלים האלו שהן בעצם נוסחת נסיגה עבור סדרת וקטורי הזיכרון m_t.' `
    'אולם הדאטה שלנו דיסקרטי (טוקנים נגיד) אז צריך לעשות דיסקרטיזציה (דגימה) לגישה הזו. כלומר במקום פונקציות רציפות תהיה לנו סדרת הקלט u_t וקטור הזיכרון m_t. גם מטריצות במערכת הדינמית שלנו צריכות לעבור דיסקרטיזציה (השערוך הרגיל של הנגזרת/גרדיאנט) ואז נקבל נוסחה רקורסיבית עבור m_t כפונקציה של u_t ו- m_t-1.  ניתן לתאר את את הדגימות עד t=T על ידי נוסחת נסיגה הזו.'

# --- 6) Delete the long stretch of old body paragraphs (old paragraphs
#        10 through 21 inclusive), keeping the blank paragraph 9 and the
#        final summary paragraph (old 22) intact. ---------------------------
$delStart = $d.Paragraphs.Item(10).Range.Start
$delEnd = $d.Paragraphs.Item(22).Range.Start
$d.Range($delStart, $delEnd).Delete()

# --- 7) The final (summary) paragraph is now Paragraphs.Item(10) -----------
Replace-ParagraphText $d.Paragraphs.Item(10) `
    'אוקיי, עכשיו סיכום במשפט אחד של המאמר הדי כבד הזה. המחברים בנו מסגרת מתמטית למידול בעיית הזיכרון של פונקציית קלט שישמש אותנו מאחורי הקלעים לבניית מודלי attention כל הדרך לממבה. ' `
    'זהו זה - יש לנו רשת בסגנון RNN כאשר הזיכרון ממודל על ידי דיסקרטיזציה של מערכת דינמית, המחשבת מקדמים של פולינומי Legendre ובאופן זה עבד לא רע אי שם ב 2020.'

# --- 8) Insert the four new paragraphs right after the secondary link
#        (paragraph 4): three blank "Normal" paragraphs followed by the new
#        introductory Hebrew paragraph. Doing this last keeps every index
#        used above stable while the edits above were applied. -------------
$d.Paragraphs.Item(4).Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs.Item(5).Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs.Item(6).Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs.Item(7).Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs.Item(8).Range.Text = 'המאמר הראשון בסדרה שלנו מנסה לטפל בעיה הראשונה של RNNs כלומר באי יכולת של רשתות אלו לדחוס את הזיכרון (קלט בחלון ההקשר) בצורה מספיק טובה. המאמר מציע גישה מקורית ומעניינת שמקורה במערכות דינמיות (Dynamic Systems) לבניית ייצוג הזיכרון. נניח שיש לנו פונקציית קלט רציפה (u(t ואנו רוצים לבנות מערכת ש״זוכרת את הפונקציה זו״ כלומר בונה ייצוג כך שיהיה אפשר לשחזרה באופן מדויק. תזכרו שכדי לתאר קלט דיסקרטי כמו טקסט אנו צריכים רק לעשות דיסקרטיזציה או לדגום את הפונקציה הזו.'

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
